$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.09425133333333334
$ws.Range("H2").Value = 0.282754
$ws.Range("I2").Value = 0.02715992817009031
$ws.Range("J2").Value = 0.02715992817009031
$ws.Range("M2").Value = 74.98881166666666
$ws.Range("N2").Value = 224.966435
$ws.Range("O2").Value = 0.6650357743745379
$ws.Range("P2").Value = 0.6650357743745379
$ws.Range("Q2").Value = 7.067795484665556
$ws.Range("R2").Value = 63.61015936199
$ws.Range("S2").Value = 0.01806232386255283
$ws.Range("T2").Value = 0.01806232386255283
$ws.Range("G3").Value = 0.09425133333333334
$ws.Range("H3").Value = 0.282754
$ws.Range("I3").Value = 0.02715992817009031
$ws.Range("J3").Value = 0.02715992817009031
$ws.Range("M3").Value = 9.661646333333335
$ws.Range("O3").Value = 0.08568398816056159
$ws.Range("P3").Value = 0.08568398816056158
$ws.Range("Q3").Value = 0.910623049111778
$ws.Range("R3").Value = 8.195607442006001
$ws.Range("S3").Value = 0.002327170963767721
$ws.Range("T3").Value = 0.002327170963767721
$ws.Range("G4").Value = 0.09425133333333334
$ws.Range("H4").Value = 0.282754
$ws.Range("I4").Value = 0.02715992817009031
$ws.Range("J4").Value = 0.02715992817009031
$ws.Range("M4").Value = 28.10860633333333
$ws.Range("N4").Value = 84.325819
$ws.Range("O4").Value = 0.2492802374649006
$ws.Range("P4").Value = 0.2492802374649006
$ws.Range("Q4").Value = 2.649273625058445
$ws.Range("R4").Value = 23.843462625526
$ws.Range("S4").Value = 0.006770433343769754
$ws.Range("T4").Value = 0.006770433343769753
$ws.Range("I5").Value = 0.357039508851706
$ws.Range("J5").Value = 0.357039508851706
$ws.Range("M5").Value = 74.98881166666666
$ws.Range("N5").Value = 224.966435
$ws.Range("O5").Value = 0.6650357743745379
$ws.Range("P5").Value = 0.6650357743745379
$ws.Range("Q5").Value = 92.91196253192832
$ws.Range("R5").Value = 836.2076627873549
$ws.Range("S5").Value = 0.237444046251499
$ws.Range("T5").Value = 0.237444046251499
$ws.Range("I6").Value = 0.357039508851706
$ws.Range("J6").Value = 0.357039508851706
$ws.Range("M6").Value = 9.661646333333335
$ws.Range("O6").Value = 0.08568398816056159
$ws.Range("P6").Value = 0.08568398816056158
$ws.Range("S6").Value = 0.0305925690493023
$ws.Range("T6").Value = 0.0305925690493023
$ws.Range("I7").Value = 0.357039508851706
$ws.Range("J7").Value = 0.357039508851706
$ws.Range("M7").Value = 28.10860633333333
$ws.Range("N7").Value = 84.325819
$ws.Range("O7").Value = 0.2492802374649006
$ws.Range("P7").Value = 0.2492802374649006
$ws.Range("Q7").Value = 34.82687244166966
$ws.Range("S7").Value = 0.08900289355090472
$ws.Range("T7").Value = 0.08900289355090472
$ws.Range("I8").Value = 0.6158005629782037
$ws.Range("J8").Value = 0.6158005629782037
$ws.Range("M8").Value = 74.98881166666666
$ws.Range("N8").Value = 224.966435
$ws.Range("O8").Value = 0.6650357743745379
$ws.Range("P8").Value = 0.6650357743745379
$ws.Range("Q8").Value = 160.2490408374811
$ws.Range("R8").Value = 1442.24136753733
$ws.Range("S8").Value = 0.4095294042604861
$ws.Range("T8").Value = 0.4095294042604861
$ws.Range("I9").Value = 0.6158005629782037
$ws.Range("J9").Value = 0.6158005629782037
$ws.Range("M9").Value = 9.661646333333335
$ws.Range("O9").Value = 0.08568398816056159
$ws.Range("P9").Value = 0.08568398816056158
$ws.Range("S9").Value = 0.05276424814749157
$ws.Range("T9").Value = 0.05276424814749157
$ws.Range("I10").Value = 0.6158005629782037
$ws.Range("J10").Value = 0.6158005629782037
$ws.Range("M10").Value = 28.10860633333333
$ws.Range("N10").Value = 84.325819
$ws.Range("O10").Value = 0.2492802374649006
$ws.Range("P10").Value = 0.2492802374649006
$ws.Range("Q10").Value = 60.06732343242689
$ws.Range("R10").Value = 540.605910891842
$ws.Range("S10").Value = 0.1535069105702261
$ws.Range("T10").Value = 0.1535069105702261
